# Apply the "Updated cryptos list" refresh (Wed Apr 12 15:31:32 UTC 2023, GitHub Actions).
# Only the cells whose values actually changed are touched; column D is forced to
# Text format first so numeric-looking price strings (e.g. "8.950") keep their exact
# text representation instead of being auto-coerced to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.973.72"
$ws.Range("E2").Value = "  -0.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.919.58"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.07"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5036"
$ws.Range("E7").Value = "  -2.87%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4028"
$ws.Range("E8").Value = "  -0.95%  "

$ws.Range("E9").Value = "  -3.10%  "

$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.04"
$ws.Range("E11").Value = "  -1.85%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.62"
$ws.Range("E12").Value = "  +0.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.912.78"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.404"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.292"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.11"
$ws.Range("E17").Value = "  -3.39%  "

$ws.Range("E18").Value = "  -1.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06483"
$ws.Range("E19").Value = "  -3.18%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.33"
$ws.Range("E20").Value = "  -1.01%  "

$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.953"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.026.78"
$ws.Range("E23").Value = "  -0.66%  "

$ws.Range("E25").Value = "  -1.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.22"
$ws.Range("E26").Value = "  +4.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.134.34"
$ws.Range("E27").Value = "  -0.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.24"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.343"
$ws.Range("E29").Value = "  -3.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.97"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.126"
$ws.Range("E31").Value = "  +2.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1042"
$ws.Range("E32").Value = "  -2.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.982"
$ws.Range("E33").Value = "  -0.61%  "

$ws.Range("E34").Value = "  +4.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02445"
$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.397"
$ws.Range("E36").Value = "  +4.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06447"
$ws.Range("E37").Value = "  -2.02%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.950"
$ws.Range("E38").Value = "  +1.39%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2161"
$ws.Range("E39").Value = "  -2.42%  "

$ws.Range("E40").Value = "  -3.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6423"
$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.38"
$ws.Range("E42").Value = "  -4.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.217"
$ws.Range("E43").Value = "  -1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.27"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5999"
$ws.Range("E46").Value = "  -2.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.163"
$ws.Range("E47").Value = "  +3.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.640"
$ws.Range("E48").Value = "  -2.67%  "

$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.215"
$ws.Range("E49").Value = "  -2.56%  "

$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.80"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.89"
$ws.Range("E51").Value = "  -0.87%  "
